$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new relationship row for RETAIL-STORE-4 supplied by WAREHOUSE-C
$ws.Range("A22").Value = "dtmi:isa95:space:CustomerLocation;1"
$ws.Range("B22").Value = "RETAIL-STORE-4"
$ws.Range("C22").Value = "WAREHOUSE-C"
$ws.Range("D22").Value = "supplies"

# Update the active selection to match the edited cell
$ws.Range("C22").Select()
